$wb = $excel.ActiveWorkbook

# Target column width (OOXML "width" units) for the Status columns, taken
# from the canonical diff: 17.2159881591797. Excel only stores ColumnWidth
# quantized to whole pixels, so we feed in the ColumnWidth value that lands
# on the closest achievable stored width (17.1666...), using the midpoint of
# the quantization bucket so we're robust to rounding.
$statusColumnWidth = 16.333333

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-06 05:15:18"
$overview.Columns.Item(5).ColumnWidth = $statusColumnWidth
$overview.Columns.Item(6).ColumnWidth = $statusColumnWidth

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-06 05:15:10"
$zhcn.Columns.Item(3).ColumnWidth = $statusColumnWidth

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-06 05:15:18"
$dede.Columns.Item(3).ColumnWidth = $statusColumnWidth

$wb.Save()
